$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row insertions for new Task Enumeration content (rows 36-51) ---
# Original layout: row37 = "Bring Up" (A37), row38 = "Test " (A38); row 36 unused/blank.

# Step 1: Insert 3 rows before "Bring Up" (row 37) -> new rows 37,38,39; "Bring Up" moves to row 40
$ws.Range("A37:A39").EntireRow.Insert()

# Step 2: Insert 5 rows before "Test " (now at row 41) -> new rows 41-45; "Test " moves to row 46
$ws.Range("A41:A45").EntireRow.Insert()

# --- Fill in the new cell values ---
$ws.Range("A36").Value = "Assembly"
$ws.Range("B37").Value = "Board prep"
$ws.Range("B38").Value = "Component Soldering"
$ws.Range("B39").Value = "Rework"

$ws.Range("B41").Value = "Set Up"
$ws.Range("B42").Value = "Base Line Existing Hardware"
$ws.Range("B43").Value = "Board Induction"
$ws.Range("B44").Value = "Smoke Test"
$ws.Range("B45").Value = "Debug Cycle"

$ws.Range("B47").Value = "Test Setup"
$ws.Range("B48").Value = "Testing the Test"
$ws.Range("B49").Value = "Testing "
$ws.Range("B50").Value = "Analysis Development "
$ws.Range("B51").Value = "Analysis"

# --- Column B width (bestfit widened to fit the new longer labels) ---
$ws.Columns.Item(2).ColumnWidth = 25.1640625

# --- Row height tweak for header rows 1-4 (16.5 -> 16) ---
$ws.Rows.Item(1).RowHeight = 16
$ws.Rows.Item(2).RowHeight = 16
$ws.Rows.Item(3).RowHeight = 16
$ws.Rows.Item(4).RowHeight = 16

# --- Sheet view / selection updates ---
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("A44").Select() | Out-Null

# --- Workbook window size/position ---
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 0
$win.Width = 28800
$win.Height = 17460
